# Review_140.docx -> Review_139.docx content swap
# (title/link swapped to a different paper, and the Hebrew review body
#  replaced with the write-up for the new paper; the old write-up is
#  cleared out to an empty run, matching the target diff.)

$d = $word.ActiveDocument

# --- 1. Heading paragraph: title text + huggingface/arxiv link -------------
$found = $d.Content.Find.Execute(
    "Review 140: [Short] One Wide Feedforward is All You Need, 07.09.2023",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Review 139: Unsupervised Compositional Concepts Discovery with Text-to-Image Generative Models, 06.09.2023",
    2)
if (-not $found) { throw "Could not find the Review 140 title text" }

$found = $d.Content.Find.Execute(
    "https://huggingface.co/papers/2309.01826",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "https://arxiv.org/abs/2306.05357.pdf",
    2)
if (-not $found) { throw "Could not find the huggingface papers link" }

# --- 2. Bold "Paper:" line --------------------------------------------------
$found = $d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.01826v2",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2306.05357v2",
    2)
if (-not $found) { throw "Could not find the bold Paper: link line" }

# --- 3. Previously-empty paragraph gets the new Hebrew review body ---------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:t xml:space="preserve">מודלי דיפוזיה גנרטיביים מסוגלים ליצור תמונות (ולא רק) באיכות מרהיבה בהתאם לתיאור טקסטואלי. עכשיו השאלה האם ניתן לגרום למודלי דיפוזיה לעשות את הפעולה ההפוכה (סוג של). כלומר להפיק קונספטים ויזואליים (כלב, מכונית וכאלו) מסט תמונות נתון. </w:t><w:br/><w:br/><w:t xml:space="preserve">היום ב-#shorthebrewpapereviews אנחנו סוקרים מאמר שמציע שיטה אלגנטית אינטואיטיבית לכך. קודם כל המאמר מגדיר מושג ״מודל דיפוזיה מרוכב״ (composite diffusion model או CDM). אתם זוכרים שבמהלך הגנרוט של תמונה אם מודלי דיפוזיה רגילים אנו מתחילים מרעש טהור ואז מסירים מנות קטנות של רעש, המחושב עם מודל מאומן, באיטרציות עד שמגיעים לתמונה נקיה. </w:t><w:br/><w:br/><w:t xml:space="preserve">המודל שמשערך את הרעש מקבל גם את התיאור של התמונה (ייצוגו הוקטורי). עם CDM אנו יוצרים תמונה עם יותר מקונספט ויזואלי אחד כאשר הרעש המנוקה הוא סכום של שערוכי רעש בהינתן קונספט מינוס הרעש המשוערך ללא התניה בתיאור (unconditioned). זה די דומה (מקרה פרטי) של classifier-free guidance שיטה פופולרית לגנרוט תמונות עם מודלי דיפוזיה. </w:t><w:br/><w:br/><w:t xml:space="preserve">אבל בואו נחזור לשאלה המקורית: איך ניתן להפיק קונספטים ויזואליים (או ייצוגם) מסט תמונות נתון? נניח שאנו רוצים ״ללמוד״ K קונספטים ויזואליים מסט של תמונות. אז כל תמונה ניתנת לייצוג בתור CDM כאשר במקום סכום של שערוכי הרעש בהינתן כל קונספט בונים את סכומם הקמור (סכום המקדמים שווה ל 1 וכולם אי שליליים) עם המקדמים נלמדים עבור כל תמונה (גם כאן מחסירים מהסכום את שערוך הרעש ה-unconditioned). </w:t><w:br/><w:br/><w:t>אז בהינתן סט תמונות מכיילים מודל דיפוזיה כאשר המטרה ללמוד את ייצוגי הקונספטים הויזואליים וגם המקדמים של כל קונספט בכל תמונה. הלוס מכמת עד כמה טוב הצלחנו לשחזר כל תמונה. לאחר מכן ניתן לגנרט תמונה עבור כל קונספט בנפרד או כל מיני שילובים שלהם.</w:t></w:r></w:p>')

# --- 4. Old Hebrew review body paragraph is cleared to an empty run --------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r/></w:p>')

Write-Output "Edit complete: $($d.Paragraphs.Count) paragraphs"
